# Generate Report for Handoff
#
# The handback-version checker noticed that the handback for
# 4a28ded3-b3a0-46a7-a25c-d047914385c2.md is stale, so its status moves
# from "Handed back: in sync with en-US" to "Ready for handoff" (with a
# fresh generation timestamp) and a new Error Detail message records why.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80c775d4ea90782192a06b71a87991f20f3273ca/e2e/4a28ded3-b3a0-46a7-a25c-d047914385c2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a7aa88ee422ebea6cfb52dc40382e4985053bc6/e2e/4a28ded3-b3a0-46a7-a25c-d047914385c2.md."

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-03 00:52:22"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-03 00:52:17"
$zhcn.Range("P3").Value = $errorDetail
# Widen the Error Detail column so the message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-03 00:52:22"
$dede.Range("P3").Value = $errorDetail
# Widen the Error Detail column so the message is readable.
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
